$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (pushes everything else down by 2)
$ws.Rows("2:3").Insert()

# Row 2: EB 97.3 (ZA7888)
$ws.Range("A2").Value = "ZA7888"
$ws.Range("B2").Value = "'97.3"
$ws.Range("C2").Value = "April-May 2022"
$ws.Range("D2").Value = "European Parliament Spring Survey, Sport and Physical Activity, and Key Challenges of our Times - The EU in 2022"

# Row 3: EB 97.2 (ZA7887)
$ws.Range("A3").Value = "ZA7887"
$ws.Range("B3").Value = "'97.2"
$ws.Range("C3").Value = "March-April 2022"
$ws.Range("D3").Value = "Corruption and Attitudes of Europeans towards Air Quality"

# Match the final selection state recorded in the workbook
$ws.Range("D3").Select() | Out-Null
